$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Insert a new column before column B.
#    This shifts the old "Kerk's laptop" (B), "Kerk's home" (C) and
#    "Kerk's desktop" (D) data one column to the right (-> C, D, E).
$ws.Columns("B").Insert()

# 2. New header row. Enter "Daryl's computer" first so the shared-string
#    table grows in the same order as the original author's edit
#    (Daryl's computer, then Relative).
$ws.Range("F1").Value = "Daryl's computer"
$ws.Range("B1").Value = "Relative"

# 3. Baseline row additions: B2 (relative baseline) and F2 (Daryl's computer baseline)
$ws.Range("B2").Value = 1
$ws.Range("F2").Value = 20.8700066891408

# 4. Row 3 (ILAsolveLIN): new "Kerk's laptop" timing + relative formula
$ws.Range("C3").Value = 0.01024248500471
$ws.Range("B3").Formula = "=C3/C2"

# 5. Row 5 (ILAsolveGSSA): new "Daryl's computer" timing + relative formula
$ws.Range("F5").Value = 62.465973368613
$ws.Range("B5").Formula = "=F5/F2"

# 6. Row 6 (ILAsimLIN): new "Kerk's laptop" timing + relative formula
$ws.Range("C6").Value = 98346.1276661065
$ws.Range("B6").Formula = "=C6/C2"

# 7. Row 7 (ILAsimVFI): formatted placeholder cell (no value yet)
$ws.Range("B7").NumberFormat = "0.000"

# 8. Row 8 (ILAsimGSSA): new "Daryl's computer" timing + relative formula
$ws.Range("F8").Value = 12422.3521871921
$ws.Range("B8").Formula = "=F8/F2"

# 9. Row 9 (OLGsolveLIN): the old raw value (shifted into C9 by the column
#    insert) is discarded and replaced by a fresh "Kerk's home" measurement
#    in D9, with a relative formula in B9.
$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = 0.130760968874471
$ws.Range("B9").Formula = "=D9/D2"

# 10. Row 12 (OLGsimVFI): the pre-existing value (shifted into C12 by the
#     column insert) actually belongs in "Kerk's home" (D12); move it there
#     and add the relative formula in B12.
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = 2163.1343990180098
$ws.Range("B12").Formula = "=D12/D2"

# 11. Number formats for the populated column-B cells.
$ws.Range("B1").NumberFormat = "0.0000000"
$ws.Range("B2").NumberFormat = "0.0000000"
$ws.Range("B3").NumberFormat = "0.0000000"
$ws.Range("B9").NumberFormat = "0.0000000"
$ws.Range("B5").NumberFormat = "0.000"
$ws.Range("B6").NumberFormat = "0.000"
$ws.Range("B8").NumberFormat = "0.000"
$ws.Range("B12").NumberFormat = "0.000"

# 12. Update the active selection to D14, matching the author's last click.
$ws.Range("D14").Select()

$wb.Save()
